# Reading row and creating array out of it
# Update the monthly figures for the "laila" row (and friends) and make
# sure a trailing (currently blank) row is present under the data, the
# same way Excel keeps a formatted-but-empty row when the sheet grows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values read while building the array from the row
$ws.Range("B2").Value = 11
$ws.Range("B4").Value = 10
$ws.Range("B5").Value = 1

# Keep row 6 present (matches the default row height) even though it has
# no cell content yet.
$ws.Rows.Item(6).RowHeight = 14.25
